$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2..538).
# The whole column was bumped by one day: 45171 (2023-09-02) -> 45172 (2023-09-03).
$ws.Range("C2:C538").Value = 45172
